# Progress-Report.xlsx update:
#  - User Interface: % COMPLETE for "Login / Register" raised to 75%
#  - Back end: two rows ("Registration Validation" / "Fix layout, match
#    design with other pages") marked 100% complete (DONE flips via its
#    table formula)
#  - Test Cases: US #18 % COMPLETE raised to 50% and given a REMARKS note
#  - Paperworks tab left as the active sheet/tab
#  - cursor/selection position nudged on each sheet, matching the saved
#    Excel session state

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "User Interface" ---
$ws1 = $wb.Worksheets.Item("User Interface")
$ws1.Range("C4").Value = 0.75
$ws1.Range("D13").Select() | Out-Null

# --- Sheet 2: "Back end" ---
$ws2 = $wb.Worksheets.Item("Back end")
$ws2.Range("E17").Value = 1
$ws2.Range("E21").Value = 1
$ws2.Range("E15").Select() | Out-Null
$win2 = $excel.ActiveWindow
$win2.ScrollRow = 3
$win2.ScrollColumn = 1

# --- Sheet 3: "Test Cases" ---
$ws3 = $wb.Worksheets.Item("Test Cases")
$ws3.Range("D20").Value = 0.5
$ws3.Range("F20").Value = "Still to test"
$ws3.Range("F21").Select() | Out-Null

# --- Sheet 4: "Paperworks" (ends as the active sheet / active tab) ---
$ws4 = $wb.Worksheets.Item("Paperworks")
$ws4.Activate() | Out-Null
